$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H4: "32.0m 6.95s" -> "32m 6.95s"
$ws.Range("H4").Value = "32m 6.95s"

# Fill in row 6 values
$ws.Range("F6").Value = 77
$ws.Range("G6").Value = 225
$ws.Range("H6").Value = "50m 12.97s"

# Update selection to M5
$ws.Range("M5").Select()
